$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 9.125343698179966
